$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins, Losses, Ties in AD1:AF1.
# Copy the formatting from the existing last header cell (AC1 - the
# "Unnamed: 28" column) so the new header cells pick up the same bold /
# bordered / centered style used by the rest of row 1, then set the text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-56: team record values (same W/L/T for every player row,
# since the whole roster shares one team record).
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 78   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 84   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
